$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells are formatted as Text so values like dates/percentages
# are stored as literal strings (matching inlineStr in the sheet XML).
$ws.Range("D2:G2").NumberFormat = "@"
$ws.Range("D3:G3").NumberFormat = "@"
$ws.Range("D4:G4").NumberFormat = "@"
$ws.Range("D5:G5").NumberFormat = "@"
$ws.Range("D6:G6").NumberFormat = "@"
$ws.Range("D7:G7").NumberFormat = "@"
$ws.Range("D8:G8").NumberFormat = "@"
$ws.Range("D9:G9").NumberFormat = "@"
$ws.Range("D10:G10").NumberFormat = "@"
$ws.Range("D11:G11").NumberFormat = "@"
$ws.Range("D12:G12").NumberFormat = "@"
$ws.Range("D13:G13").NumberFormat = "@"
$ws.Range("D14:G14").NumberFormat = "@"
$ws.Range("D15:G15").NumberFormat = "@"
$ws.Range("D16:G16").NumberFormat = "@"
$ws.Range("D17:G17").NumberFormat = "@"
$ws.Range("D18:G18").NumberFormat = "@"
$ws.Range("D19:G19").NumberFormat = "@"
$ws.Range("D20:G20").NumberFormat = "@"
$ws.Range("D21:G21").NumberFormat = "@"
$ws.Range("D22:G22").NumberFormat = "@"
$ws.Range("D23:G23").NumberFormat = "@"
$ws.Range("D24:G24").NumberFormat = "@"
$ws.Range("D25:G25").NumberFormat = "@"
$ws.Range("D26:G26").NumberFormat = "@"
$ws.Range("D27:G27").NumberFormat = "@"
$ws.Range("F28:G28").NumberFormat = "@"
$ws.Range("F29:G29").NumberFormat = "@"
$ws.Range("F30:G30").NumberFormat = "@"
$ws.Range("F31:G31").NumberFormat = "@"
$ws.Range("F32:G32").NumberFormat = "@"
$ws.Range("F33:G33").NumberFormat = "@"
$ws.Range("F34:G34").NumberFormat = "@"
$ws.Range("F35:G35").NumberFormat = "@"
$ws.Range("F36:G36").NumberFormat = "@"
$ws.Range("F37:G37").NumberFormat = "@"
$ws.Range("F38:G38").NumberFormat = "@"
$ws.Range("F39:G39").NumberFormat = "@"
$ws.Range("D40:G40").NumberFormat = "@"
$ws.Range("D41:G41").NumberFormat = "@"
$ws.Range("D42:G42").NumberFormat = "@"
$ws.Range("D43:G43").NumberFormat = "@"
$ws.Range("D44:G44").NumberFormat = "@"
$ws.Range("D45:G45").NumberFormat = "@"
$ws.Range("D46:G46").NumberFormat = "@"
$ws.Range("D47:G47").NumberFormat = "@"
$ws.Range("F48:G48").NumberFormat = "@"
$ws.Range("F49:G49").NumberFormat = "@"
$ws.Range("F50:G50").NumberFormat = "@"
$ws.Range("F51:G51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "275.19"
$ws.Range("E2").Value = "-0.27%"
$ws.Range("F2").Value = "10-1-2023"
$ws.Range("G2").Value = "0"
# Row 3
$ws.Range("D3").Value = "26.92"
$ws.Range("E3").Value = "-0.80%"
$ws.Range("F3").Value = "10-1-2023"
$ws.Range("G3").Value = "0"
# Row 4
$ws.Range("D4").Value = "4.907"
$ws.Range("E4").Value = "2.94%"
$ws.Range("F4").Value = "10-1-2023"
$ws.Range("G4").Value = "0"
# Row 5
$ws.Range("D5").Value = "0.06334"
$ws.Range("E5").Value = "2.01%"
$ws.Range("F5").Value = "10-1-2023"
$ws.Range("G5").Value = "0"
# Row 6
$ws.Range("D6").Value = "6.915"
$ws.Range("E6").Value = "1.71%"
$ws.Range("F6").Value = "10-1-2023"
$ws.Range("G6").Value = "0"
# Row 7
$ws.Range("D7").Value = "3.349"
$ws.Range("E7").Value = "3.65%"
$ws.Range("F7").Value = "10-1-2023"
$ws.Range("G7").Value = "0"
# Row 8
$ws.Range("D8").Value = "1.286"
$ws.Range("E8").Value = "39.38%"
$ws.Range("F8").Value = "10-1-2023"
$ws.Range("G8").Value = "0"
# Row 9
$ws.Range("D9").Value = "0.8821"
$ws.Range("E9").Value = "1.83%"
$ws.Range("F9").Value = "10-1-2023"
$ws.Range("G9").Value = "0"
# Row 10
$ws.Range("D10").Value = "0.1471"
$ws.Range("E10").Value = "1.40%"
$ws.Range("F10").Value = "10-1-2023"
$ws.Range("G10").Value = "0"
# Row 11
$ws.Range("D11").Value = "0.05086"
$ws.Range("E11").Value = "-2.76%"
$ws.Range("F11").Value = "10-1-2023"
$ws.Range("G11").Value = "0"
# Row 12
$ws.Range("D12").Value = "0.07361"
$ws.Range("E12").Value = "1.31%"
$ws.Range("F12").Value = "10-1-2023"
$ws.Range("G12").Value = "0"
# Row 13
$ws.Range("D13").Value = "0.03144"
$ws.Range("E13").Value = "0.45%"
$ws.Range("F13").Value = "10-1-2023"
$ws.Range("G13").Value = "0"
# Row 14
$ws.Range("D14").Value = "0.09025"
$ws.Range("E14").Value = "-0.32%"
$ws.Range("F14").Value = "10-1-2023"
$ws.Range("G14").Value = "0"
# Row 15
$ws.Range("D15").Value = "0.001555"
$ws.Range("E15").Value = "0.31%"
$ws.Range("F15").Value = "10-1-2023"
$ws.Range("G15").Value = "0"
# Row 16
$ws.Range("D16").Value = "0.0006300"
$ws.Range("E16").Value = "1.96%"
$ws.Range("F16").Value = "10-1-2023"
$ws.Range("G16").Value = "0"
# Row 17
$ws.Range("D17").Value = "0.006036"
$ws.Range("E17").Value = "-0.30%"
$ws.Range("F17").Value = "10-1-2023"
$ws.Range("G17").Value = "0"
# Row 18
$ws.Range("D18").Value = "3.457"
$ws.Range("E18").Value = "-0.53%"
$ws.Range("F18").Value = "10-1-2023"
$ws.Range("G18").Value = "0"
# Row 19
$ws.Range("D19").Value = "2.282"
$ws.Range("E19").Value = "0.33%"
$ws.Range("F19").Value = "10-1-2023"
$ws.Range("G19").Value = "0"
# Row 20
$ws.Range("D20").Value = "0.3165"
$ws.Range("E20").Value = "2.57%"
$ws.Range("F20").Value = "10-1-2023"
$ws.Range("G20").Value = "0"
# Row 21
$ws.Range("D21").Value = "0.1313"
$ws.Range("E21").Value = "0.48%"
$ws.Range("F21").Value = "10-1-2023"
$ws.Range("G21").Value = "0"
# Row 22
$ws.Range("D22").Value = "3.911"
$ws.Range("E22").Value = "1.76%"
$ws.Range("F22").Value = "10-1-2023"
$ws.Range("G22").Value = "0"
# Row 23
$ws.Range("D23").Value = "0.04362"
$ws.Range("E23").Value = "2.66%"
$ws.Range("F23").Value = "10-1-2023"
$ws.Range("G23").Value = "0"
# Row 24
$ws.Range("D24").Value = "0.001176"
$ws.Range("E24").Value = "0.14%"
$ws.Range("F24").Value = "10-1-2023"
$ws.Range("G24").Value = "0"
# Row 25
$ws.Range("D25").Value = "0.003687"
$ws.Range("E25").Value = "-12.46%"
$ws.Range("F25").Value = "10-1-2023"
$ws.Range("G25").Value = "0"
# Row 26
$ws.Range("D26").Value = "0.0001199"
$ws.Range("F26").Value = "10-1-2023"
$ws.Range("G26").Value = "0"
# Row 27
$ws.Range("D27").Value = "0.0001694"
$ws.Range("E27").Value = "-12.36%"
$ws.Range("F27").Value = "10-1-2023"
$ws.Range("G27").Value = "0"
# Row 28
$ws.Range("F28").Value = "10-1-2023"
$ws.Range("G28").Value = "0"
# Row 29
$ws.Range("F29").Value = "10-1-2023"
$ws.Range("G29").Value = "0"
# Row 30
$ws.Range("F30").Value = "10-1-2023"
$ws.Range("G30").Value = "0"
# Row 31
$ws.Range("F31").Value = "10-1-2023"
$ws.Range("G31").Value = "0"
# Row 32
$ws.Range("F32").Value = "10-1-2023"
$ws.Range("G32").Value = "0"
# Row 33
$ws.Range("F33").Value = "10-1-2023"
$ws.Range("G33").Value = "0"
# Row 34
$ws.Range("F34").Value = "10-1-2023"
$ws.Range("G34").Value = "0"
# Row 35
$ws.Range("F35").Value = "10-1-2023"
$ws.Range("G35").Value = "0"
# Row 36
$ws.Range("F36").Value = "10-1-2023"
$ws.Range("G36").Value = "0"
# Row 37
$ws.Range("F37").Value = "10-1-2023"
$ws.Range("G37").Value = "0"
# Row 38
$ws.Range("F38").Value = "10-1-2023"
$ws.Range("G38").Value = "0"
# Row 39
$ws.Range("F39").Value = "10-1-2023"
$ws.Range("G39").Value = "0"
# Row 40
$ws.Range("D40").Value = "0.04050"
$ws.Range("E40").Value = "1.24%"
$ws.Range("F40").Value = "10-1-2023"
$ws.Range("G40").Value = "0"
# Row 41
$ws.Range("D41").Value = "0.006608"
$ws.Range("E41").Value = "6.58%"
$ws.Range("F41").Value = "10-1-2023"
$ws.Range("G41").Value = "0"
# Row 42
$ws.Range("D42").Value = "0.1165"
$ws.Range("E42").Value = "2.76%"
$ws.Range("F42").Value = "10-1-2023"
$ws.Range("G42").Value = "0"
# Row 43
$ws.Range("D43").Value = "0.002188"
$ws.Range("E43").Value = "2.78%"
$ws.Range("F43").Value = "10-1-2023"
$ws.Range("G43").Value = "0"
# Row 44
$ws.Range("D44").Value = "0.01267"
$ws.Range("E44").Value = "6.68%"
$ws.Range("F44").Value = "10-1-2023"
$ws.Range("G44").Value = "0"
# Row 45
$ws.Range("D45").Value = "0.00005300"
$ws.Range("E45").Value = "3.58%"
$ws.Range("F45").Value = "10-1-2023"
$ws.Range("G45").Value = "0"
# Row 46
$ws.Range("D46").Value = "2.356"
$ws.Range("E46").Value = "163.31%"
$ws.Range("F46").Value = "10-1-2023"
$ws.Range("G46").Value = "0"
# Row 47
$ws.Range("D47").Value = "0.02595"
$ws.Range("E47").Value = "5.08%"
$ws.Range("F47").Value = "10-1-2023"
$ws.Range("G47").Value = "0"
# Row 48
$ws.Range("F48").Value = "10-1-2023"
$ws.Range("G48").Value = "0"
# Row 49
$ws.Range("F49").Value = "10-1-2023"
$ws.Range("G49").Value = "0"
# Row 50
$ws.Range("F50").Value = "10-1-2023"
$ws.Range("G50").Value = "0"
# Row 51
$ws.Range("F51").Value = "10-1-2023"
$ws.Range("G51").Value = "0"
